$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Test autorisatieprofielen" to "Mapping API-GBA"
# (the workbook is being repurposed as a feature for filling address
# lines / mapping the API response to GBA fields).
$ws.Name = "Mapping API-GBA"

# Move the cursor to where the author left off editing (F65), which
# updates the active cell / selection recorded for the frozen
# bottom-right pane.
$ws.Range("F65").Select()
